$wb = $excel.ActiveWorkbook

# ===== sheet1 (index 1) =====
$ws = $wb.Worksheets.Item(1)
$ws.Range("B1").Value = "name"
$ws.Range("C1").Value = "area"
$ws.Range("D1").Value = "share_portion"
$ws.Range("E1").Value = "owner"
$ws.Range("F1").Value = "register_date"
$ws.Range("G1").Value = "register_reason"
$ws.Range("H1").Value = "acquire_value"
$ws.Range("I1").Value = "property_category"
$ws.Range("J1").Value = "category"
$ws.Range("K1").Value = "date"
$ws.Range("L1").Value = "legislator_name"
$ws.Range("M1").Value = "legislator_id"
$ws.Range("N1").Value = "source_file"
$ws.Range("O1").Value = "index"

$ws.Range("A2").Value = 13
$ws.Range("B2").Value = "臺南市安南區海東段00450069地號"
$ws.Range("C2").Value = 85.59
$ws.Range("D2").Value = "全部"
$ws.Range("E2").Value = "許添財"
$ws.Range("F2").Value = "87年07月03日"
$ws.Range("G2").Value = "買賣"
$ws.Range("H2").Value = "(超過五年交付第一銀行）"
$ws.Range("I2").Value = "land"
$ws.Range("J2").Value = "normal"
$ws.Range("K2").NumberFormat = "@"
$ws.Range("K2").Value = "2012-03-22"
$ws.Range("L2").Value = "許添財"
$ws.Range("M2").Value = 639
$ws.Range("N2").Value = "tmpb8d31"
$ws.Range("O2").Value = 13
$ws.Range("A3").Value = 14
$ws.Range("B3").Value = "臺南市安南區海東段00450070地號"
$ws.Range("C3").Value = 85.59
$ws.Range("D3").Value = "全部"
$ws.Range("E3").Value = "洪淑貞"
$ws.Range("F3").Value = "88年08月16日"
$ws.Range("G3").Value = "買賣"
$ws.Range("H3").Value = "(超過12年）"
$ws.Range("I3").Value = "land"
$ws.Range("J3").Value = "normal"
$ws.Range("K3").NumberFormat = "@"
$ws.Range("K3").Value = "2012-03-22"
$ws.Range("L3").Value = "許添財"
$ws.Range("M3").Value = 639
$ws.Range("N3").Value = "tmpb8d31"
$ws.Range("O3").Value = 14
$ws.Range("A4").Value = 15
$ws.Range("B4").Value = "165thSt.FlushingNewYorkBlock5250Lot26"
$ws.Range("C4").Value = 464.5
$ws.Range("D4").Value = "全部"
$ws.Range("E4").Value = "許添財洪淑貞"
$ws.Range("G4").Value = "買賣"
$ws.Range("H4").Value = "(超過20年）"
$ws.Range("I4").Value = "land"
$ws.Range("J4").Value = "normal"
$ws.Range("K4").NumberFormat = "@"
$ws.Range("K4").Value = "2012-03-22"
$ws.Range("L4").Value = "許添財"
$ws.Range("M4").Value = 639
$ws.Range("N4").Value = "tmpb8d31"
$ws.Range("O4").Value = 15

$ws.Range("H1").Copy()
$ws.Range("I1:O1").PasteSpecial(-4122)
$ws.Range("H2").Copy()
$ws.Range("I2:O2").PasteSpecial(-4122)
$ws.Range("H3").Copy()
$ws.Range("I3:O3").PasteSpecial(-4122)
$ws.Range("H4").Copy()
$ws.Range("I4:O4").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# ===== sheet2 (index 2) =====
$ws = $wb.Worksheets.Item(2)
$ws.Range("B1").Value = "建物標示"
$ws.Range("C1").Value = "面積（平方公尺）"
$ws.Range("D1").Value = "權利範圍(持分）"
$ws.Range("E1").Value = "所有權人"
$ws.Range("F1").Value = "登記（取得）時間"
$ws.Range("G1").Value = "登記（取得）原乱"
$ws.Range("H1").Value = "取得價額"

$ws.Range("A2").Value = 20
$ws.Range("B2").Value = "臺南市安南區海東段01774000建號"
$ws.Range("C2").Value = 282.77
$ws.Range("D2").Value = "全部"
$ws.Range("E2").Value = "許添財"
$ws.Range("F2").Value = "88年間"
$ws.Range("G2").Value = "買賣"
$ws.Range("H2").Value = "(超過五年交付第一銀行含陽台）"
$ws.Range("A3").Value = 21
$ws.Range("B3").Value = "臺南市安南區海東段01773000建號"
$ws.Range("C3").Value = 282.77
$ws.Range("D3").Value = "全部"
$ws.Range("E3").Value = "洪淑貞"
$ws.Range("F3").Value = "88年08月16日"
$ws.Range("G3").Value = "買賣"
$ws.Range("H3").Value = "(超過12年含陽台）"
$ws.Range("A4").Value = 22
$ws.Range("B4").Value = "165thSt.FlushingNewYorkBlock5250Lot26"
$ws.Range("C4").Value = 465.5
$ws.Range("D4").Value = "全部"
$ws.Range("E4").Value = "許添財洪淑貞"
$ws.Range("G4").Value = "買賣"
$ws.Range("H4").Value = "(超過20年）"

# ===== sheet3 (index 3) =====
$ws = $wb.Worksheets.Item(3)
$ws.Range("B1").Value = "廠牌型號"
$ws.Range("C1").Value = "汽缸容量"
$ws.Range("D1").Value = "所有人"
$ws.Range("E1").Value = "登記（取得）時間"
$ws.Range("F1").Value = "登記（取得）原因"
$ws.Range("G1").Value = "取得價額"

$ws.Range("A2").Value = 32
$ws.Range("B2").Value = "國瑞MVIEPE"
$ws.Range("C2").Value = 2995
$ws.Range("D2").Value = "許添財"
$ws.Range("E2").Value = "93年09月23日"
$ws.Range("F2").Value = "買賣"
$ws.Range("G2").Value = 1100000

# ===== sheet4 (index 4) =====
$ws = $wb.Worksheets.Item(4)
$ws.Range("B1").Value = "存放機構(應敘明分支機構）"
$ws.Range("C1").Value = "種類"
$ws.Range("D1").Value = "幣別"
$ws.Range("E1").Value = "所有人"
$ws.Range("F1").Value = "外幣總額"
$ws.Range("G1").Value = "新臺幣總額或折合新臺幣總額"

$ws.Range("A2").Value = 48
$ws.Range("B2").Value = "臺灣中小企業銀行安平分行"
$ws.Range("C2").Value = "活期儲蓄存款"
$ws.Range("D2").Value = "新臺幣"
$ws.Range("E2").Value = "許添財"
$ws.Range("G2").Value = 2165459
$ws.Range("A3").Value = 49
$ws.Range("B3").Value = "臺灣銀行群賢分行"
$ws.Range("C3").Value = "綜合存款"
$ws.Range("D3").Value = "新臺幣"
$ws.Range("E3").Value = "洪淑頁"
$ws.Range("G3").Value = 2253656
$ws.Range("A4").Value = 50
$ws.Range("B4").Value = "臺灣銀行南都分行"
$ws.Range("C4").Value = "活期儲蓄存款"
$ws.Range("D4").Value = "新臺幣"
$ws.Range("E4").Value = "許添財"
$ws.Range("G4").Value = 1779236
$ws.Range("A5").Value = 51
$ws.Range("B5").Value = "臺灣銀行南都分行"
$ws.Range("C5").Value = "定期存款"
$ws.Range("D5").Value = "新臺幣"
$ws.Range("E5").Value = "許添財"
$ws.Range("G5").Value = 3000000
$ws.Range("A6").Value = 52
$ws.Range("B6").Value = "臺灣銀行南都分行"
$ws.Range("C6").Value = "綜合存款"
$ws.Range("D6").Value = "新臺幣"
$ws.Range("E6").Value = "洪淑貞"
$ws.Range("G6").Value = 1000
$ws.Range("A7").Value = 53
$ws.Range("B7").Value = "台北富邦商業銀行駐立分行"
$ws.Range("C7").Value = "活期儲蓄存款"
$ws.Range("D7").Value = "新臺幣"
$ws.Range("E7").Value = "許添財"
$ws.Range("G7").Value = 254
$ws.Range("A8").Value = 54
$ws.Range("B8").Value = "華南商業銀行台南分行"
$ws.Range("C8").Value = "活期存款"
$ws.Range("D8").Value = "新臺幣"
$ws.Range("E8").Value = "許添財"
$ws.Range("G8").Value = 1210
$ws.Range("A9").Value = 55
$ws.Range("B9").Value = "華南商業銀行東台南分行"
$ws.Range("C9").Value = "活期存款"
$ws.Range("D9").Value = "新臺幣"
$ws.Range("E9").Value = "許添財"
$ws.Range("G9").Value = 1845
$ws.Range("A10").Value = 56
$ws.Range("B10").Value = "第一商業銀行運河分行"
$ws.Range("C10").Value = "活期儲蓄存款"
$ws.Range("D10").Value = "新臺幣"
$ws.Range("E10").Value = "許添財"
$ws.Range("G10").Value = 11567
$ws.Range("A11").Value = 57
$ws.Range("B11").Value = "合作金庫商業銀行南興分行"
$ws.Range("C11").Value = "支票存款"
$ws.Range("D11").Value = "新臺幣"
$ws.Range("E11").Value = "許添財"
$ws.Range("G11").Value = 1977
$ws.Range("A12").Value = 58
$ws.Range("B12").Value = "金城商業銀行府城分行"
$ws.Range("C12").Value = "活期存款"
$ws.Range("D12").Value = "新臺幣"
$ws.Range("E12").Value = "許添財"
$ws.Range("G12").Value = 388
$ws.Range("A13").Value = 59
$ws.Range("B13").Value = "臺灣銀行(註2)"
$ws.Range("C13").Value = "綜合存款"
$ws.Range("D13").Value = "美金"
$ws.Range("E13").Value = "許添財"
$ws.Range("F13").Value = 369.68
$ws.Range("G13").Value = 10991
$ws.Range("A14").Value = 60
$ws.Range("B14").Value = "美商花旗銀行"
$ws.Range("C14").Value = "支票存款"
$ws.Range("D14").Value = "美金"
$ws.Range("E14").Value = "許添財洪淑貞"
$ws.Range("F14").Value = 2000.38
$ws.Range("G14").Value = 59475
$ws.Range("A15").Value = 61
$ws.Range("B15").Value = "日商三菱東京曰聯銀行"
$ws.Range("C15").Value = "其他存款"
$ws.Range("D15").Value = "美金"
$ws.Range("E15").Value = "洪淑貞"
$ws.Range("F15").Value = 100954.8
$ws.Range("G15").Value = 3001588
$ws.Range("A16").Value = 62
$ws.Range("B16").Value = "美商花旗銀行"
$ws.Range("C16").Value = "其他存款"
$ws.Range("D16").Value = "美金"
$ws.Range("E16").Value = "許添財"
$ws.Range("F16").Value = 16557.19
$ws.Range("G16").Value = 492278
$ws.Range("A17").Value = 63
$ws.Range("B17").Value = "美商花旗銀行"
$ws.Range("C17").Value = "其他存款"
$ws.Range("D17").Value = "美金"
$ws.Range("E17").Value = "洪淑貞"
$ws.Range("F17").Value = 4889.73
$ws.Range("G17").Value = 145381
$ws.Range("A18").Value = 64
$ws.Range("B18").Value = "美商摩根大通銀行"
$ws.Range("C18").Value = "其他存款"
$ws.Range("D18").Value = "美金"
$ws.Range("E18").Value = "許添財"
$ws.Range("F18").Value = 4201.16
$ws.Range("G18").Value = 124909
$ws.Range("A19").Value = 65
$ws.Range("B19").Value = "美商摩根大通銀行"
$ws.Range("C19").Value = "其他存款"
$ws.Range("D19").Value = "美金"
$ws.Range("E19").Value = "洪淑貞"
$ws.Range("F19").Value = 4201.16
$ws.Range("G19").Value = 124909

# ===== sheet5 (index 5) =====
$ws = $wb.Worksheets.Item(5)
$ws.Range("B1").Value = "名"
$ws.Range("C1").Value = "稱"
$ws.Range("D1").Value = "所"
$ws.Range("E1").Value = "有"
$ws.Range("F1").Value = "人"
$ws.Range("G1").Value = "單位數"

$ws.Range("A2").Value = 89
$ws.Range("B2").Value = "(九）珠寶古董字畫及另"
$ws.Range("C2").Value = "他具有相當價值之財產（總fl"
$ws.Range("E2").Value = "貢額:新臺幣1657元）"
$ws.Range("A3").Value = 90
$ws.Range("B3").Value = "財產種類"
$ws.Range("C3").Value = "項"
$ws.Range("D3").Value = "件"
$ws.Range("E3").Value = "所有人"
$ws.Range("F3").Value = "價"
$ws.Range("G3").Value = "額"
$ws.Range("A4").Value = 91
$ws.Range("B4").Value = "黃金1克(註3)"
$ws.Range("C4").Value = 1
$ws.Range("E4").Value = "許添財"
$ws.Range("F4").Value = 1657.1
